# Update the "想去人数" (column F) figures on the "展览" and "全部类型"
# sheets to reflect the latest scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# New values for column F, keyed by row number, per worksheet.
# (Row 21 differed slightly between the two sheets before the update, so
# each sheet gets its own explicit mapping to be safe.)

$updatesExhibition = @{
    2  = 313
    3  = 13896
    5  = 100
    6  = 185
    7  = 283
    8  = 498
    10 = 90
    13 = 53
    14 = 459
    15 = 5893
    16 = 139
    17 = 93
    18 = 983
    19 = 122
    20 = 62
    21 = 159
    22 = 283
}

$updatesAllTypes = @{
    2  = 313
    3  = 13897
    5  = 100
    6  = 185
    7  = 283
    8  = 498
    10 = 90
    13 = 53
    14 = 459
    15 = 5893
    16 = 139
    17 = 93
    18 = 983
    19 = 122
    20 = 62
    21 = 159
    22 = 283
}

$sheetUpdates = @{
    "展览"   = $updatesExhibition
    "全部类型" = $updatesAllTypes
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $updates = $sheetUpdates[$sheetName]
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
